$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

$ws.Range("C$row").NumberFormat = "@"
$ws.Range("E$row").NumberFormat = "@"
$ws.Range("F$row").NumberFormat = "@"

$ws.Range("A$row").Value = "gadar"
$ws.Range("B$row").Value = "ACtion"
$ws.Range("C$row").Value = "3"
$ws.Range("D$row").Value = "sunny deol"
$ws.Range("E$row").Value = "50"
$ws.Range("F$row").Value = "3"
